$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set each changed cell explicitly, forcing text storage
# for columns that contain numeric-looking strings (D) so Excel does not
# reinterpret them as numbers and strip formatting (e.g. "1.00" -> 1).

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.488.59'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.90%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.227.73'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.32%  '

# Row 4
$ws.Range('E4').Value = '  +0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '111.51'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -7.01%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '293.71'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +9.92%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.621'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.36%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.00%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.47'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -8.15%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0919'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.05%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.22'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.49%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.86'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.20%  '

# Row 14
$ws.Range('E14').Value = '  +11.81%  '

# Row 15
$ws.Range('E15').Value = '  -2.86%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.11'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.79%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.561.06'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.45%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.229.20'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.35%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.469.87'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.92%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.37'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +6.04%  '

# Row 21
$ws.Range('E21').Value = '  -3.93%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.07'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.03%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.47'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +20.12%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.53%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '229.77'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.69%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.68%  '

# Row 27
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.57%  '

# Row 28
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.66'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.07%  '

# Row 29
$ws.Range('E29').Value = '  -0.88%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.48'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -10.55%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.27'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.60%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.53'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.13%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.09'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.93%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0902'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.62%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.23'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +13.18%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.71'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.99%  '

# Row 37
$ws.Range('E37').Value = '  +1.47%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.127'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.36%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0377'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.72%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.104'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.77%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.42'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.78%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.53'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.18%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.235'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.78%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.05%  '

# Row 45
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.75'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.04%  '

# Row 46
$ws.Range('E46').Value = '  -5.17%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.46'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.75%  '

# Row 48
$ws.Range('E48').Value = '  +3.68%  '

# Row 49
$ws.Range('E49').Value = '  +0.77%  '

# Row 50
$ws.Range('E50').Value = '  +0.91%  '

# Row 51
$ws.Range('E51').Value = '  +7.29%  '
